# Commit: "all canine test cases 72"
#
# Adds three new worksheets to the workbook, mirroring the "stat" query
# output that the Commons automation tool appends alongside the existing
# CypherOutput / Message sheets:
#   - CypherOutput_Message : a copy of the "Message" sheet (connection +
#                             cypher-query log) for the CypherOutput run.
#   - StatOutput            : the numeric summary returned by the stats
#                              cypher query (file/sample/case/study counts).
#   - StatOutput_Message    : the connection + cypher-query log for BOTH the
#                              original CypherOutput query and the new stats
#                              query, one after another.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Pull the existing "Message" sheet's values so we don't have to retype
# the (long, quote-heavy) Neo4j connection info / cypher query text.
# ---------------------------------------------------------------------
$msgSheet = $wb.Worksheets.Item("Message")

$neo4jUrlLabel   = $msgSheet.Cells.Item(1, 1).Value()
$neo4jUrlValue   = $msgSheet.Cells.Item(2, 1).Value()
$userLabel       = $msgSheet.Cells.Item(3, 1).Value()
$userValue       = $msgSheet.Cells.Item(4, 1).Value()
$pwdLabel        = $msgSheet.Cells.Item(5, 1).Value()
$pwdValue        = $msgSheet.Cells.Item(6, 1).Value()
$cypherLabel     = $msgSheet.Cells.Item(7, 1).Value()
$cypherQuery     = $msgSheet.Cells.Item(8, 1).Value()
$outputLabel     = $msgSheet.Cells.Item(9, 1).Value()
$outputPath      = $msgSheet.Cells.Item(10, 1).Value()

# The new stats cypher query (same shape as the original, but returns
# aggregate counts instead of per-case rows).
$statQuery = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.sex IN [''Female Phenotype'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

function Add-SheetAtEnd([string]$name) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $s = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
    $s.Name = $name
    return $s
}

function Write-MessageBlock($sheet, [int]$startRow) {
    $sheet.Cells.Item($startRow + 0, 1).Value = $neo4jUrlLabel
    $sheet.Cells.Item($startRow + 1, 1).Value = $neo4jUrlValue
    $sheet.Cells.Item($startRow + 2, 1).Value = $userLabel
    $sheet.Cells.Item($startRow + 3, 1).Value = $userValue
    $sheet.Cells.Item($startRow + 4, 1).Value = $pwdLabel
    $sheet.Cells.Item($startRow + 5, 1).Value = $pwdValue
    $sheet.Cells.Item($startRow + 6, 1).Value = $cypherLabel
}

# ---------------------------------------------------------------------
# 1) CypherOutput_Message -- exact copy of the Message sheet.
# ---------------------------------------------------------------------
$cypherOutputMessage = Add-SheetAtEnd "CypherOutput_Message"
Write-MessageBlock $cypherOutputMessage 1
$cypherOutputMessage.Cells.Item(8, 1).Value = $cypherQuery
$cypherOutputMessage.Cells.Item(9, 1).Value = $outputLabel
$cypherOutputMessage.Cells.Item(10, 1).Value = $outputPath

# ---------------------------------------------------------------------
# 2) StatOutput -- header row of stat field names + one row of counts.
#    The counts are numeric-looking text (not numbers) in the source
#    workbook, so force text formatting before assigning them.
# ---------------------------------------------------------------------
$statOutput = Add-SheetAtEnd "StatOutput"
$statOutput.Cells.Item(1, 1).Value = "number_of_files"
$statOutput.Cells.Item(1, 2).Value = "number_of_sample"
$statOutput.Cells.Item(1, 3).Value = "number_of_cases"
$statOutput.Cells.Item(1, 4).Value = "number_of_study"

$statCells = @(
    $statOutput.Cells.Item(2, 1),
    $statOutput.Cells.Item(2, 2),
    $statOutput.Cells.Item(2, 3),
    $statOutput.Cells.Item(2, 4)
)
$statValues = @("0", "0", "5", "1")
for ($i = 0; $i -lt $statCells.Length; $i++) {
    $statCells[$i].NumberFormat = "@"
    $statCells[$i].Value = $statValues[$i]
}

# ---------------------------------------------------------------------
# 3) StatOutput_Message -- the CypherOutput message block followed by a
#    second block for the stats query.
# ---------------------------------------------------------------------
$statOutputMessage = Add-SheetAtEnd "StatOutput_Message"
Write-MessageBlock $statOutputMessage 1
$statOutputMessage.Cells.Item(8, 1).Value = $cypherQuery
$statOutputMessage.Cells.Item(9, 1).Value = $outputLabel
$statOutputMessage.Cells.Item(10, 1).Value = $outputPath

Write-MessageBlock $statOutputMessage 11
$statOutputMessage.Cells.Item(18, 1).Value = $statQuery
$statOutputMessage.Cells.Item(19, 1).Value = $outputLabel
$statOutputMessage.Cells.Item(20, 1).Value = $outputPath

# ---------------------------------------------------------------------
# Restore the original active sheet/tab selection (adding sheets moves
# selection to the newly created sheet).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("CypherOutput").Activate()
